$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 0. Remove the existing "_GoBack" bookmark (it will be re-added later
#    at the location of the final edit made in this session, which is
#    how Word tracks the last-edit position).
# ---------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------
# 1. "Connect Data to PID" + "# (2)" -> single run "Connect Data to PID# (2)"
#    (the bookmark that used to separate them is already gone above, so
#    we just need to force the two runs to coalesce into one by doing a
#    no-op replace across the whole phrase)
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Connect Data to PID# (2)", $false, $false, $false, $false, $false, $true, 1, $false, "Connect Data to PID# (2)", 2) | Out-Null

# ---------------------------------------------------------------------
# 2. ".    As a Developer, I want a barebones template for the Android
#    app as a starting." -> "...starting place." and re-establish the
#    "_GoBack" bookmark at the very end of that paragraph.
# ---------------------------------------------------------------------

# 2a. Force a run split between "template fo" and "r the Android..." by
#     temporarily bookmarking that point before editing the text, so the
#     edit below does not coalesce that boundary away.
$findSplit = $d.Content
$findSplit.Find.Execute(".    As a Developer, I want a barebones template fo") | Out-Null
$splitPoint = $d.Range($findSplit.End, $findSplit.End)
$d.Bookmarks.Add("TempSplit", $splitPoint) | Out-Null

# 2b. Replace "starting." with "starting place."
$findText = $d.Content
$findText.Find.Execute("starting.") | Out-Null
$findText.Text = "starting place."

# 2c. Remove the temporary split-marker bookmark (the run boundary it
#     created remains intact).
$d.Bookmarks("TempSplit").Delete()

# 2d. Find the (now updated) end of that paragraph's text.
$findEnd = $d.Content
$findEnd.Find.Execute("starting place.") | Out-Null
$endPos = $findEnd.End

# 2e. Work around a boundary issue when bookmarking exactly at a
#     paragraph's trailing (pilcrow) position: temporarily extend the
#     paragraph by one placeholder character, bookmark before it, then
#     remove the placeholder again.
$placeholder = $d.Range($endPos, $endPos)
$placeholder.InsertAfter("Z")
$bmRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
$d.Range($endPos, $endPos + 1).Delete()
